# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# YDS sheet: per-play rush/pass yardage logs (offense + defense).
# Week 16 (and the simulated remainder of the season) appended new
# play-by-play yardage figures to each running list.
# -----------------------------------------------------------------
$ydsWs = $wb.Sheets("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 4 5 1 1 1 3 1 3 -4 1 -3 1 3 0 -4 5"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 5 2 7 1 40 4 14 6 8 20 5 13 4 8"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 3 3 7 4 6 5 3 3 -2 0 -4 8 7 3 8 19 4 8 5 5 4 8 -4 3 0 4 4 4 10 6 0 4 13 3 5"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 8 5 16 18 6 4 0 15 10 20 -5 5 7 13 12 12 18 6 3 28"

# -----------------------------------------------------------------
# OFF sheet: down/distance play counts + passing/sack/4th-down totals
# -----------------------------------------------------------------
$offWs = $wb.Sheets("OFF")

$offWs.Range("B2").Value = 6
$offWs.Range("C2").Value = 182
$offWs.Range("D2").Value = 11
$offWs.Range("E2").Value = 17
$offWs.Range("F2").Value = 61
$offWs.Range("G2").Value = 55
$offWs.Range("I2").Value = 8
$offWs.Range("J2").Value = 26
$offWs.Range("N2").Value = 16
$offWs.Range("O2").Value = 24
$offWs.Range("P2").Value = 15

$offWs.Range("C3").Value = 154
$offWs.Range("E3").Value = 33
$offWs.Range("F3").Value = 87
$offWs.Range("H3").Value = 35
$offWs.Range("I3").Value = 66
$offWs.Range("L3").Value = 233
$offWs.Range("M3").Value = 157
$offWs.Range("Q3").Value = 430

# -----------------------------------------------------------------
# DEF sheet: same shape as OFF, opponent-facing totals
# -----------------------------------------------------------------
$defWs = $wb.Sheets("DEF")

$defWs.Range("C2").Value = 174
$defWs.Range("D2").Value = 5
$defWs.Range("F2").Value = 63
$defWs.Range("G2").Value = 48
$defWs.Range("J2").Value = 33
$defWs.Range("N2").Value = 10
$defWs.Range("O2").Value = 24
$defWs.Range("P2").Value = 9

$defWs.Range("C3").Value = 142
$defWs.Range("F3").Value = 78
$defWs.Range("G3").Value = 38
$defWs.Range("H3").Value = 27
$defWs.Range("I3").Value = 47
$defWs.Range("J3").Value = 58
$defWs.Range("L3").Value = 220
$defWs.Range("M3").Value = 126
$defWs.Range("Q3").Value = 411

# -----------------------------------------------------------------
# ST sheet: special-teams totals + per-kick logs (KO/PT/FG distance
# buckets)
# -----------------------------------------------------------------
$stWs = $wb.Sheets("ST")

$stWs.Range("B2").Value = 69
$stWs.Range("D2").Value = 59
$stWs.Range("F2").Value = 255
$stWs.Range("G2").Value = 248
$stWs.Range("L2").Value = 61
$stWs.Range("M2").Value = 51
$stWs.Range("N2").Value = 58
$stWs.Range("O2").Value = 31
$stWs.Range("B3").Value = 57

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 62"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 23"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 18"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 54 28 56 40 45"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 0 6 0 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0"

# -----------------------------------------------------------------
# TURNS sheet: turnover totals
# -----------------------------------------------------------------
$turnsWs = $wb.Sheets("TURNS")

$turnsWs.Range("C3").Value = 5
$turnsWs.Range("D3").Value = 7
$turnsWs.Range("E3").Value = 7

# -----------------------------------------------------------------
# PEN sheet: penalty counts
# -----------------------------------------------------------------
$penWs = $wb.Sheets("PEN")

$penWs.Range("B2").Value = 15
